# Update crypto symbol list snapshot (prices + a few volume-label tweaks).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain numeric-looking text (t="inlineStr").
# Force Text format first so Excel doesn't silently convert these back into
# floating point numbers when we re-assign the string values.
$ws.Range("D2:D25").NumberFormat = "@"
$ws.Range("D40:D50").NumberFormat = "@"

$ws.Range("D2").Value  = "247.18"
$ws.Range("D3").Value  = "22.63"
$ws.Range("D4").Value  = "5.291"
$ws.Range("D5").Value  = "0.05721"
$ws.Range("D6").Value  = "3.426"
$ws.Range("D7").Value  = "0.8054"
$ws.Range("D8").Value  = "0.8635"
$ws.Range("D9").Value  = "0.1422"
$ws.Range("D10").Value = "0.07345"
$ws.Range("D11").Value = "0.03043"
$ws.Range("D12").Value = "0.03114"
$ws.Range("D13").Value = "0.09390"
$ws.Range("D14").Value = "3.882"
$ws.Range("D15").Value = "0.001589"
$ws.Range("D16").Value = "0.04820"
$ws.Range("D17").Value = "0.0005852"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "0.006156"
$ws.Range("D20").Value = "0.0009965"
$ws.Range("D22").Value = "3.707"
$ws.Range("D23").Value = "6.316"
$ws.Range("D24").Value = "2.195"
$ws.Range("D25").Value = "0.3275"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "0.03938"
$ws.Range("D41").Value = "0.006770"
$ws.Range("D43").Value = "0.003200"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").Value = "0.007974"
$ws.Range("D45").Value = "0.00005612"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "0.1819"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.01010"
